$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph (it directly follows
#    the "Play Gods of Asgard Megaways free: Review & Details"
#    Heading1 paragraph at the top of the document).
# ------------------------------------------------------------------
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Meta description")) {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -gt 0) {
    $d.Paragraphs.Item($metaIndex).Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new paragraph ("Play Gods of Asgard Megaways free:
#    Review & Details", bold) right before the final paragraph in the
#    document (the italic image-generation-prompt paragraph).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

# A second, empty paragraph is appended to the fragment purely so the
# XML importer splits it off from the following (italic) paragraph
# instead of merging the runs into it; the stray paragraph left behind
# by this is removed right afterwards.
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Gods of Asgard Megaways free: Review &amp; Details</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/></w:p>'
$insPoint.InsertXML($xmlFrag)

$strayPara = $d.Paragraphs.Item($count + 1)
$strayPara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the old image-generation prompt text (now in the final,
#    italic paragraph) with the meta-description copy that used to sit
#    at the top of the document.
# ------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute(
    "Create an eye-catching feature image for Gods of Asgard Megaways! The image should be in cartoon style and prominently feature a happy Maya warrior with glasses. Make use of vibrant colors and catchy graphics to draw in potential players. Additionally, include elements from the Norse mythology theme, such as lightning bolts, wolf symbols, and ravens. The image should convey the excitement and epic nature of the game while also showcasing its unique features, like the Megaways format and the God symbols that offer exciting bonuses.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Gods of Asgard Megaways slot game. Play it for free, and enjoy its stunning graphics, epic soundtrack, and various bonus features!",
    2)
